$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.181.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.60%  '

$ws.Range("D3").Value = "'1.685.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'216.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = "'22.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.13%  '

$ws.Range("E9").Value = '  +2.59%  '

$ws.Range("D10").Value = "'0.0625"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.27%  '

$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").Value = "'1.924.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = "'1.686.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("E14").Value = '  +2.18%  '

$ws.Range("D15").Value = "'0.558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.32%  '

$ws.Range("D16").Value = "'66.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.88%  '

$ws.Range("D17").Value = "'27.195.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").Value = "'235.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '

$ws.Range("D19").Value = "'8.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.45%  '

$ws.Range("D20").Value = "'0.0₃0743"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.70%  '

$ws.Range("D22").Value = "'4.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.19%  '

$ws.Range("D23").Value = "'9.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.39%  '

$ws.Range("E24").Value = '  -2.19%  '

$ws.Range("D25").Value = "'146.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").Value = "'7.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.84%  '

$ws.Range("E27").Value = '  -1.63%  '

$ws.Range("E28").Value = '  +0.25%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").Value = "'0.0504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.06%  '

$ws.Range("D31").Value = "'1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").Value = "'3.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.90%  '

$ws.Range("D33").Value = "'1.543.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("E34").Value = '  +1.72%  '

$ws.Range("D35").Value = "'1.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.98%  '

$ws.Range("E36").Value = '  +2.31%  '

$ws.Range("D37").Value = "'0.946"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.01%  '

$ws.Range("D38").Value = "'2.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.51%  '

$ws.Range("E39").Value = '  -1.28%  '

$ws.Range("E40").Value = '  +2.22%  '

$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("D42").Value = "'69.12"
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("E44").Value = '  -0.66%  '

$ws.Range("D45").Value = "'1.831.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("D46").Value = "'0.790"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.16%  '

$ws.Range("D47").Value = "'90.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("E48").Value = '  +8.83%  '

$ws.Range("D49").Value = "'1.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.77%  '

$ws.Range("D50").Value = "'8.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.33%  '

$ws.Range("E51").Value = '  -0.81%  '
